# Actualización automática 2025-10-25 10:30:09
#
# Updates sales figures for client "PADILLA MIER BERTHA MARIETA" (PORCELANATO,
# october) and "TULCAN NARVAEZ EDITH MARITZA" (240X80 PORCELANATO, INODOROS,
# PIEDRA SINTERIZADA, october) across the three report sheets, plus the
# "N de 21" completion counters and the compliance-summary totals that are
# derived from them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO" — per-client sales by product group
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# PADILLA MIER BERTHA MARIETA — PORCELANATO
$wsGrupo.Range("M17").Value = 9155.459999999999

# TULCAN NARVAEZ EDITH MARITZA — 240X80 PORCELANATO / INODOROS / PIEDRA SINTERIZADA
$wsGrupo.Range("D22").Value = 1526.4
$wsGrupo.Range("H22").Value = 1940.3
$wsGrupo.Range("L22").Value = 1611.96

# Row 23 "N de 21" completion counters (one per product-group column)
$wsGrupo.Range("D23").Value = "2 de 21"
$wsGrupo.Range("H23").Value = "1 de 21"
$wsGrupo.Range("L23").Value = "1 de 21"
$wsGrupo.Range("M23").Value = "3 de 21"

# ---------------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL" — per-client sales by month
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# PADILLA MIER BERTHA MARIETA — octubre
$wsMensual.Range("F17").Value = 9155.459999999999

# TULCAN NARVAEZ EDITH MARITZA — octubre
$wsMensual.Range("F22").Value = 5078.66

# Totals row — octubre
$wsMensual.Range("F23").Value = 17753.64

# ---------------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL" — budget compliance summary by product group
# ---------------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column widths widened by one character to fit the larger figures
$wsCumpl.Columns.Item(4).ColumnWidth = 13.17
$wsCumpl.Columns.Item(5).ColumnWidth = 22.17

# Row 3 — 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 2079.36
$wsCumpl.Range("E3").Value = 3425.25890386263
$wsCumpl.Range("F3").Value = 0.3777482213238956

# Row 6 — INODOROS
$wsCumpl.Range("D6").Value = 1940.3
$wsCumpl.Range("E6").Value = 967.2836814602599
$wsCumpl.Range("F6").Value = 0.6673238718362643

# Row 11 — PIEDRA SINTERIZADA
$wsCumpl.Range("D11").Value = 1611.96
$wsCumpl.Range("E11").Value = 4232.48916370549
$wsCumpl.Range("F11").Value = 0.2758104236769487

# Row 12 — PORCELANATO
$wsCumpl.Range("D12").Value = 12122.02
$wsCumpl.Range("E12").Value = 25617.72
$wsCumpl.Range("F12").Value = 0.3212004110256192

# Row 14 — Total
$wsCumpl.Range("D14").Value = 17753.64
$wsCumpl.Range("E14").Value = 37671.1014788039
$wsCumpl.Range("F14").Value = 0.3203197620107896
